$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = "participant_id"

$ws.Range("B3").Select()
